# Scheduled-runner market data refresh for the Leve-profit tables.
# Updates currentAveragePrice / NQ / HQ price & profit columns (H:N) on
# each profession sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# the latest pulled market data, row by row.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 821.5454999999999  # H18
$ws.Cells.Item(18, 9).Value = 813.7  # I18
$ws.Cells.Item(18, 10).Value = 900  # J18
$ws.Cells.Item(18, 11).Value = 813.7  # K18
$ws.Cells.Item(18, 12).Value = 900  # L18
$ws.Cells.Item(18, 13).Value = -529.7  # M18
$ws.Cells.Item(18, 14).Value = -1468  # N18
$ws.Cells.Item(40, 8).Value = 1792.5883  # H40
$ws.Cells.Item(40, 10).Value = 2073.3333  # J40
$ws.Cells.Item(40, 12).Value = 2073.3333  # L40
$ws.Cells.Item(40, 14).Value = -2423.3333  # N40
$ws.Cells.Item(129, 8).Value = 738.12  # H129
$ws.Cells.Item(129, 10).Value = 1054.5454  # J129
$ws.Cells.Item(129, 12).Value = 3163.6362  # L129
$ws.Cells.Item(129, 14).Value = -13163.6362  # N129
$ws.Cells.Item(138, 8).Value = 6900222.5  # H138
$ws.Cells.Item(138, 9).Value = 2330.8462  # I138
$ws.Cells.Item(138, 10).Value = 12504760  # J138
$ws.Cells.Item(138, 11).Value = 6992.5386  # K138
$ws.Cells.Item(138, 12).Value = 37514280  # L138
$ws.Cells.Item(138, 13).Value = -1852.5386  # M138
$ws.Cells.Item(138, 14).Value = -37524560  # N138
$ws.Cells.Item(141, 8).Value = 2081.7083  # H141
$ws.Cells.Item(141, 9).Value = 1993.7391  # I141
$ws.Cells.Item(141, 10).Value = 4105  # J141
$ws.Cells.Item(141, 11).Value = 5981.2173  # K141
$ws.Cells.Item(141, 12).Value = 12315  # L141
$ws.Cells.Item(141, 13).Value = -801.2173000000003  # M141
$ws.Cells.Item(141, 14).Value = -22675  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2110.75  # H45
$ws.Cells.Item(45, 9).Value = 1860.6666  # I45
$ws.Cells.Item(45, 10).Value = 2260.8  # J45
$ws.Cells.Item(45, 11).Value = 1860.6666  # K45
$ws.Cells.Item(45, 12).Value = 2260.8  # L45
$ws.Cells.Item(45, 13).Value = -1483.6666  # M45
$ws.Cells.Item(45, 14).Value = -3014.8  # N45
$ws.Cells.Item(74, 8).Value = 144714.86  # H74
$ws.Cells.Item(74, 9).Value = 251278  # I74
$ws.Cells.Item(74, 10).Value = 2630.6667  # J74
$ws.Cells.Item(74, 11).Value = 251278  # K74
$ws.Cells.Item(74, 12).Value = 2630.6667  # L74
$ws.Cells.Item(74, 13).Value = -250404  # M74
$ws.Cells.Item(74, 14).Value = -4378.6667  # N74
$ws.Cells.Item(77, 8).Value = 144714.86  # H77
$ws.Cells.Item(77, 9).Value = 251278  # I77
$ws.Cells.Item(77, 10).Value = 2630.6667  # J77
$ws.Cells.Item(77, 11).Value = 1256390  # K77
$ws.Cells.Item(77, 12).Value = 13153.3335  # L77
$ws.Cells.Item(77, 13).Value = -1252022  # M77
$ws.Cells.Item(77, 14).Value = -21889.3335  # N77
$ws.Cells.Item(112, 8).Value = 20819.428  # H112
$ws.Cells.Item(112, 10).Value = 20819.428  # J112
$ws.Cells.Item(112, 12).Value = 20819.428  # L112
$ws.Cells.Item(112, 14).Value = -23773.428  # N112
$ws.Cells.Item(122, 8).Value = 1910.7858  # H122
$ws.Cells.Item(122, 9).Value = 1580.3334  # I122
$ws.Cells.Item(122, 10).Value = 2505.6  # J122
$ws.Cells.Item(122, 11).Value = 4741.0002  # K122
$ws.Cells.Item(122, 12).Value = 7516.799999999999  # L122
$ws.Cells.Item(122, 13).Value = -2291.0002  # M122
$ws.Cells.Item(122, 14).Value = -12416.8  # N122
$ws.Cells.Item(125, 8).Value = 34781.668  # H125
$ws.Cells.Item(125, 10).Value = 34781.668  # J125
$ws.Cells.Item(125, 12).Value = 34781.668  # L125
$ws.Cells.Item(125, 14).Value = -44621.668  # N125
$ws.Cells.Item(133, 8).Value = 100000  # H133
$ws.Cells.Item(133, 10).Value = 100000  # J133
$ws.Cells.Item(133, 12).Value = 100000  # L133
$ws.Cells.Item(133, 14).Value = -105060  # N133

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 46493.668  # H59
$ws.Cells.Item(59, 10).Value = 46493.668  # J59
$ws.Cells.Item(59, 12).Value = 46493.668  # L59
$ws.Cells.Item(59, 14).Value = -48187.668  # N59
$ws.Cells.Item(94, 8).Value = 817.55554  # H94
$ws.Cells.Item(94, 9).Value = 773.6  # I94
$ws.Cells.Item(94, 10).Value = 872.5  # J94
$ws.Cells.Item(94, 11).Value = 773.6  # K94
$ws.Cells.Item(94, 12).Value = 872.5  # L94
$ws.Cells.Item(94, 13).Value = -322.6  # M94
$ws.Cells.Item(94, 14).Value = -1774.5  # N94
$ws.Cells.Item(134, 8).Value = 7586.933  # H134
$ws.Cells.Item(134, 9).Value = 7215.6924  # I134
$ws.Cells.Item(134, 10).Value = 10000  # J134
$ws.Cells.Item(134, 11).Value = 21647.0772  # K134
$ws.Cells.Item(134, 12).Value = 30000  # L134
$ws.Cells.Item(134, 13).Value = -19112.0772  # M134
$ws.Cells.Item(134, 14).Value = -35070  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 87.30768999999999  # H7
$ws.Cells.Item(7, 9).Value = 28.333334  # I7
$ws.Cells.Item(7, 10).Value = 105  # J7
$ws.Cells.Item(7, 11).Value = 28.333334  # K7
$ws.Cells.Item(7, 12).Value = 105  # L7
$ws.Cells.Item(7, 13).Value = 84.66666599999999  # M7
$ws.Cells.Item(7, 14).Value = -331  # N7
$ws.Cells.Item(58, 8).Value = 2347.75  # H58
$ws.Cells.Item(58, 9).Value = 2127.7778  # I58
$ws.Cells.Item(58, 10).Value = 3007.6667  # J58
$ws.Cells.Item(58, 11).Value = 2127.7778  # K58
$ws.Cells.Item(58, 12).Value = 3007.6667  # L58
$ws.Cells.Item(58, 13).Value = -1924.7778  # M58
$ws.Cells.Item(58, 14).Value = -3413.6667  # N58
$ws.Cells.Item(68, 8).Value = 23333.334  # H68
$ws.Cells.Item(68, 10).Value = 23333.334  # J68
$ws.Cells.Item(68, 12).Value = 23333.334  # L68
$ws.Cells.Item(68, 14).Value = -24831.334  # N68
$ws.Cells.Item(69, 8).Value = 19297  # H69
$ws.Cells.Item(69, 9).Value = 2000  # I69
$ws.Cells.Item(69, 10).Value = 22179.834  # J69
$ws.Cells.Item(69, 11).Value = 2000  # K69
$ws.Cells.Item(69, 12).Value = 22179.834  # L69
$ws.Cells.Item(69, 13).Value = -1251  # M69
$ws.Cells.Item(69, 14).Value = -23677.834  # N69
$ws.Cells.Item(70, 8).Value = 30222.5  # H70
$ws.Cells.Item(70, 10).Value = 30222.5  # J70
$ws.Cells.Item(70, 12).Value = 30222.5  # L70
$ws.Cells.Item(70, 14).Value = -30852.5  # N70
$ws.Cells.Item(71, 8).Value = 23333.334  # H71
$ws.Cells.Item(71, 10).Value = 23333.334  # J71
$ws.Cells.Item(71, 12).Value = 70000.00199999999  # L71
$ws.Cells.Item(71, 14).Value = -77488.00199999999  # N71
$ws.Cells.Item(72, 8).Value = 19297  # H72
$ws.Cells.Item(72, 9).Value = 2000  # I72
$ws.Cells.Item(72, 10).Value = 22179.834  # J72
$ws.Cells.Item(72, 11).Value = 6000  # K72
$ws.Cells.Item(72, 12).Value = 66539.50199999999  # L72
$ws.Cells.Item(72, 13).Value = -2256  # M72
$ws.Cells.Item(72, 14).Value = -74027.50199999999  # N72
$ws.Cells.Item(73, 8).Value = 30222.5  # H73
$ws.Cells.Item(73, 10).Value = 30222.5  # J73
$ws.Cells.Item(73, 12).Value = 30222.5  # L73
$ws.Cells.Item(73, 14).Value = -32406.5  # N73
$ws.Cells.Item(74, 8).Value = 27666.584  # H74
$ws.Cells.Item(74, 10).Value = 27666.584  # J74
$ws.Cells.Item(74, 12).Value = 27666.584  # L74
$ws.Cells.Item(74, 14).Value = -29414.584  # N74
$ws.Cells.Item(75, 8).Value = 29866.666  # H75
$ws.Cells.Item(75, 10).Value = 29866.666  # J75
$ws.Cells.Item(75, 12).Value = 29866.666  # L75
$ws.Cells.Item(75, 14).Value = -31862.666  # N75
$ws.Cells.Item(77, 8).Value = 27666.584  # H77
$ws.Cells.Item(77, 10).Value = 27666.584  # J77
$ws.Cells.Item(77, 12).Value = 82999.75199999999  # L77
$ws.Cells.Item(77, 14).Value = -91735.75199999999  # N77
$ws.Cells.Item(78, 8).Value = 29866.666  # H78
$ws.Cells.Item(78, 10).Value = 29866.666  # J78
$ws.Cells.Item(78, 12).Value = 89599.99800000001  # L78
$ws.Cells.Item(78, 14).Value = -99583.99800000001  # N78
$ws.Cells.Item(81, 8).Value = 28863.334  # H81
$ws.Cells.Item(81, 10).Value = 28863.334  # J81
$ws.Cells.Item(81, 12).Value = 28863.334  # L81
$ws.Cells.Item(81, 14).Value = -30859.334  # N81
$ws.Cells.Item(82, 8).Value = 29963.334  # H82
$ws.Cells.Item(82, 10).Value = 29963.334  # J82
$ws.Cells.Item(82, 12).Value = 29963.334  # L82
$ws.Cells.Item(82, 14).Value = -30685.334  # N82
$ws.Cells.Item(84, 8).Value = 28863.334  # H84
$ws.Cells.Item(84, 10).Value = 28863.334  # J84
$ws.Cells.Item(84, 12).Value = 86590.00199999999  # L84
$ws.Cells.Item(84, 14).Value = -96574.00199999999  # N84
$ws.Cells.Item(85, 8).Value = 29963.334  # H85
$ws.Cells.Item(85, 10).Value = 29963.334  # J85
$ws.Cells.Item(85, 12).Value = 29963.334  # L85
$ws.Cells.Item(85, 14).Value = -32459.334  # N85
$ws.Cells.Item(87, 8).Value = 26998  # H87
$ws.Cells.Item(87, 10).Value = 26998  # J87
$ws.Cells.Item(87, 12).Value = 26998  # L87
$ws.Cells.Item(87, 14).Value = -29370  # N87
$ws.Cells.Item(90, 8).Value = 26998  # H90
$ws.Cells.Item(90, 10).Value = 26998  # J90
$ws.Cells.Item(90, 12).Value = 80994  # L90
$ws.Cells.Item(90, 14).Value = -92850  # N90
$ws.Cells.Item(136, 8).Value = 2347.75  # H136
$ws.Cells.Item(136, 9).Value = 2127.7778  # I136
$ws.Cells.Item(136, 10).Value = 3007.6667  # J136
$ws.Cells.Item(136, 11).Value = 6383.3334  # K136
$ws.Cells.Item(136, 12).Value = 9023.000100000001  # L136
$ws.Cells.Item(136, 13).Value = -3833.3334  # M136
$ws.Cells.Item(136, 14).Value = -14123.0001  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 859.0599999999999  # H131
$ws.Cells.Item(131, 9).Value = 467.22223  # I131
$ws.Cells.Item(131, 10).Value = 897.8132000000001  # J131
$ws.Cells.Item(131, 11).Value = 1401.66669  # K131
$ws.Cells.Item(131, 12).Value = 2693.4396  # L131
$ws.Cells.Item(131, 13).Value = 3638.33331  # M131
$ws.Cells.Item(131, 14).Value = -12773.4396  # N131
$ws.Cells.Item(132, 8).Value = 2732.7827  # H132
$ws.Cells.Item(132, 9).Value = 2492.7  # I132
$ws.Cells.Item(132, 10).Value = 4333.3335  # J132
$ws.Cells.Item(132, 11).Value = 22434.3  # K132
$ws.Cells.Item(132, 12).Value = 39000.0015  # L132
$ws.Cells.Item(132, 13).Value = -19904.3  # M132
$ws.Cells.Item(132, 14).Value = -44060.0015  # N132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 201161.4  # H122
$ws.Cells.Item(122, 9).Value = 251076.75  # I122
$ws.Cells.Item(122, 10).Value = 1500  # J122
$ws.Cells.Item(122, 11).Value = 753230.25  # K122
$ws.Cells.Item(122, 12).Value = 4500  # L122
$ws.Cells.Item(122, 13).Value = -750780.25  # M122
$ws.Cells.Item(122, 14).Value = -9400  # N122
$ws.Cells.Item(132, 8).Value = 2488.8064  # H132
$ws.Cells.Item(132, 9).Value = 2353.4375  # I132
$ws.Cells.Item(132, 10).Value = 2633.2  # J132
$ws.Cells.Item(132, 11).Value = 7060.3125  # K132
$ws.Cells.Item(132, 12).Value = 7899.599999999999  # L132
$ws.Cells.Item(132, 13).Value = -4530.3125  # M132
$ws.Cells.Item(132, 14).Value = -12959.6  # N132
$ws.Cells.Item(136, 8).Value = 31900  # H136
$ws.Cells.Item(136, 10).Value = 31900  # J136
$ws.Cells.Item(136, 12).Value = 95700  # L136
$ws.Cells.Item(136, 14).Value = -100800  # N136

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 46252.773  # H22
$ws.Cells.Item(22, 9).Value = 500250  # I22
$ws.Cells.Item(22, 10).Value = 853.05  # J22
$ws.Cells.Item(22, 11).Value = 500250  # K22
$ws.Cells.Item(22, 12).Value = 853.05  # L22
$ws.Cells.Item(22, 13).Value = -499955  # M22
$ws.Cells.Item(22, 14).Value = -1443.05  # N22
$ws.Cells.Item(27, 8).Value = 46252.773  # H27
$ws.Cells.Item(27, 9).Value = 500250  # I27
$ws.Cells.Item(27, 10).Value = 853.05  # J27
$ws.Cells.Item(27, 11).Value = 500250  # K27
$ws.Cells.Item(27, 12).Value = 853.05  # L27
$ws.Cells.Item(27, 13).Value = -500143  # M27
$ws.Cells.Item(27, 14).Value = -1067.05  # N27
$ws.Cells.Item(110, 8).Value = 20395.8  # H110
$ws.Cells.Item(110, 10).Value = 20395.8  # J110
$ws.Cells.Item(110, 12).Value = 20395.8  # L110
$ws.Cells.Item(110, 14).Value = -28575.8  # N110
$ws.Cells.Item(115, 8).Value = 22325  # H115
$ws.Cells.Item(115, 10).Value = 22325  # J115
$ws.Cells.Item(115, 12).Value = 22325  # L115
$ws.Cells.Item(115, 14).Value = -24675  # N115
$ws.Cells.Item(122, 8).Value = 3250  # H122
$ws.Cells.Item(122, 9).Value = 2925  # I122
$ws.Cells.Item(122, 10).Value = 3900  # J122
$ws.Cells.Item(122, 11).Value = 8775  # K122
$ws.Cells.Item(122, 12).Value = 11700  # L122
$ws.Cells.Item(122, 13).Value = -6325  # M122
$ws.Cells.Item(122, 14).Value = -16600  # N122
$ws.Cells.Item(132, 8).Value = 3015.8262  # H132
$ws.Cells.Item(132, 9).Value = 2549.9375  # I132
$ws.Cells.Item(132, 10).Value = 4080.7144  # J132
$ws.Cells.Item(132, 11).Value = 7649.8125  # K132
$ws.Cells.Item(132, 12).Value = 12242.1432  # L132
$ws.Cells.Item(132, 13).Value = -5119.8125  # M132
$ws.Cells.Item(132, 14).Value = -17302.1432  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 110764.37  # H122
$ws.Cells.Item(122, 9).Value = 1350  # I122
$ws.Cells.Item(122, 11).Value = 4050  # K122
$ws.Cells.Item(122, 13).Value = -1600  # M122
$ws.Cells.Item(132, 8).Value = 4194.591  # H132
$ws.Cells.Item(132, 9).Value = 4974  # I132
$ws.Cells.Item(132, 10).Value = 3259.3  # J132
$ws.Cells.Item(132, 11).Value = 14922  # K132
$ws.Cells.Item(132, 12).Value = 9777.900000000001  # L132
$ws.Cells.Item(132, 13).Value = -12392  # M132
$ws.Cells.Item(132, 14).Value = -14837.9  # N132
$ws.Cells.Item(139, 8).Value = 35000  # H139
$ws.Cells.Item(139, 9).Value = 0  # I139
$ws.Cells.Item(139, 10).Value = 35000  # J139
$ws.Cells.Item(139, 11).Value = 0  # K139
$ws.Cells.Item(139, 12).Value = 35000  # L139
$ws.Cells.Item(139, 13).ClearContents()  # M139
$ws.Cells.Item(139, 14).Value = -45280  # N139
